$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted right after the existing row 81
# (Camote, "1a (cosecha)") entry. Insert a blank row at position 82 which
# shifts the old rows 82:127 down to 83:128.
$ws.Rows("82:82").Insert()

# The row that used to be row 81 (still holding its original values,
# since Insert() only shifted the rows below it) is duplicated into the
# newly created row 82, preserving the historical record.
$ws.Range("A81:R81").Copy($ws.Range("A82"))

# Row 81 itself is then updated with the new week's figures: a later
# date and a higher reported volume.
$ws.Range("D81").Value2 = 44603
$ws.Range("J81").Value2 = 200
